$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I26").Value = "ba"
$ws.Range("J26").Value = "Appreciation"
$ws.Range("I40").Value = "b"
$ws.Range("J40").Value = "Acknowledge (Backchannel)"
$ws.Range("I43").Value = "sd"
$ws.Range("J43").Value = "Statement-non-opinion"
$ws.Range("I44").Value = "sd"
$ws.Range("J44").Value = "Statement-non-opinion"
$ws.Range("I56").Value = "ba"
$ws.Range("J56").Value = "Appreciation"
$ws.Range("I62").Value = "b"
$ws.Range("J62").Value = "Acknowledge (Backchannel)"
$ws.Range("I66").Value = "b"
$ws.Range("J66").Value = "Acknowledge (Backchannel)"
$ws.Range("I76").Value = "sv"
$ws.Range("J76").Value = "Statement-opinion"
$ws.Range("I79").Value = "ba"
$ws.Range("J79").Value = "Appreciation"
$ws.Range("I83").Value = "ba"
$ws.Range("J83").Value = "Appreciation"
$ws.Range("I98").Value = "ba"
$ws.Range("J98").Value = "Appreciation"
$ws.Range("I99").Value = "sd"
$ws.Range("J99").Value = "Statement-non-opinion"
$ws.Range("I111").Value = "ba"
$ws.Range("J111").Value = "Appreciation"
$ws.Range("I118").Value = "sv"
$ws.Range("J118").Value = "Statement-opinion"
$ws.Range("I120").Value = "ba"
$ws.Range("J120").Value = "Appreciation"
$ws.Range("I126").Value = "aa"
$ws.Range("J126").Value = "Agree/Accept"
$ws.Range("I127").Value = "sd"
$ws.Range("J127").Value = "Statement-non-opinion"
$ws.Range("I132").Value = "sv"
$ws.Range("J132").Value = "Statement-opinion"
$ws.Range("I145").Value = "b"
$ws.Range("J145").Value = "Acknowledge (Backchannel)"
$ws.Range("I167").Value = "ba"
$ws.Range("J167").Value = "Appreciation"
$ws.Range("I169").Value = "ba"
$ws.Range("J169").Value = "Appreciation"
$ws.Range("I172").Value = "ba"
$ws.Range("J172").Value = "Appreciation"
$ws.Range("I174").Value = "ba"
$ws.Range("J174").Value = "Appreciation"
$ws.Range("I176").Value = "ba"
$ws.Range("J176").Value = "Appreciation"
$ws.Range("I186").Value = "b"
$ws.Range("J186").Value = "Acknowledge (Backchannel)"
$ws.Range("I205").Value = "sd"
$ws.Range("J205").Value = "Statement-non-opinion"
$ws.Range("I212").Value = "ba"
$ws.Range("J212").Value = "Appreciation"
$ws.Range("I217").Value = "aa"
$ws.Range("J217").Value = "Agree/Accept"
$ws.Range("I237").Value = "sv"
$ws.Range("J237").Value = "Statement-opinion"
$ws.Range("I250").Value = "ba"
$ws.Range("J250").Value = "Appreciation"
$ws.Range("I281").Value = "sv"
$ws.Range("J281").Value = "Statement-opinion"
$ws.Range("I290").Value = "ba"
$ws.Range("J290").Value = "Appreciation"
$ws.Range("I292").Value = "sd"
$ws.Range("J292").Value = "Statement-non-opinion"
$ws.Range("I308").Value = "ba"
$ws.Range("J308").Value = "Appreciation"
$ws.Range("I313").Value = "sd"
$ws.Range("J313").Value = "Statement-non-opinion"
$ws.Range("I315").Value = "ba"
$ws.Range("J315").Value = "Appreciation"
$ws.Range("I316").Value = "%"
$ws.Range("J316").Value = "Uninterpretable"
$ws.Range("I317").Value = "aa"
$ws.Range("J317").Value = "Agree/Accept"
$ws.Range("I320").Value = "aa"
$ws.Range("J320").Value = "Agree/Accept"
$ws.Range("I327").Value = "b"
$ws.Range("J327").Value = "Acknowledge (Backchannel)"
$ws.Range("I336").Value = "sd"
$ws.Range("J336").Value = "Statement-non-opinion"
$ws.Range("I342").Value = "ba"
$ws.Range("J342").Value = "Appreciation"
$ws.Range("I349").Value = "aa"
$ws.Range("J349").Value = "Agree/Accept"
$ws.Range("I352").Value = "sd"
$ws.Range("J352").Value = "Statement-non-opinion"
$ws.Range("I359").Value = "sd"
$ws.Range("J359").Value = "Statement-non-opinion"
$ws.Range("I383").Value = "aa"
$ws.Range("J383").Value = "Agree/Accept"
$ws.Range("I392").Value = "sv"
$ws.Range("J392").Value = "Statement-opinion"
$ws.Range("I393").Value = "sv"
$ws.Range("J393").Value = "Statement-opinion"
$ws.Range("I396").Value = "b"
$ws.Range("J396").Value = "Acknowledge (Backchannel)"
$ws.Range("I414").Value = "sd"
$ws.Range("J414").Value = "Statement-non-opinion"
$ws.Range("I425").Value = "b"
$ws.Range("J425").Value = "Acknowledge (Backchannel)"
$ws.Range("I433").Value = "aa"
$ws.Range("J433").Value = "Agree/Accept"
$ws.Range("I447").Value = "sd"
$ws.Range("J447").Value = "Statement-non-opinion"
$ws.Range("I451").Value = "sd"
$ws.Range("J451").Value = "Statement-non-opinion"
$ws.Range("I457").Value = "b"
$ws.Range("J457").Value = "Acknowledge (Backchannel)"
$ws.Range("I458").Value = "b"
$ws.Range("J458").Value = "Acknowledge (Backchannel)"
$ws.Range("I465").Value = "sd"
$ws.Range("J465").Value = "Statement-non-opinion"
$ws.Range("I487").Value = "sd"
$ws.Range("J487").Value = "Statement-non-opinion"
$ws.Range("I504").Value = "%"
$ws.Range("J504").Value = "Uninterpretable"
$ws.Range("I505").Value = "%"
$ws.Range("J505").Value = "Uninterpretable"
$ws.Range("I512").Value = "sd"
$ws.Range("J512").Value = "Statement-non-opinion"
$ws.Range("I515").Value = "aa"
$ws.Range("J515").Value = "Agree/Accept"
$ws.Range("I516").Value = "sd"
$ws.Range("J516").Value = "Statement-non-opinion"
$ws.Range("I517").Value = "sv"
$ws.Range("J517").Value = "Statement-opinion"
$ws.Range("I526").Value = "b"
$ws.Range("J526").Value = "Acknowledge (Backchannel)"
$ws.Range("I536").Value = "ba"
$ws.Range("J536").Value = "Appreciation"
$ws.Range("I563").Value = "sv"
$ws.Range("J563").Value = "Statement-opinion"
$ws.Range("I573").Value = "sv"
$ws.Range("J573").Value = "Statement-opinion"
$ws.Range("I587").Value = "sv"
$ws.Range("J587").Value = "Statement-opinion"
$ws.Range("I594").Value = "b"
$ws.Range("J594").Value = "Acknowledge (Backchannel)"
$ws.Range("I597").Value = "b"
$ws.Range("J597").Value = "Acknowledge (Backchannel)"
$ws.Range("I602").Value = "sd"
$ws.Range("J602").Value = "Statement-non-opinion"
$ws.Range("I603").Value = "sv"
$ws.Range("J603").Value = "Statement-opinion"
$ws.Range("I616").Value = "ba"
$ws.Range("J616").Value = "Appreciation"
$ws.Range("I626").Value = "sd"
$ws.Range("J626").Value = "Statement-non-opinion"
$ws.Range("I630").Value = "sv"
$ws.Range("J630").Value = "Statement-opinion"
$ws.Range("I634").Value = "sd"
$ws.Range("J634").Value = "Statement-non-opinion"
$ws.Range("I663").Value = "sv"
$ws.Range("J663").Value = "Statement-opinion"
$ws.Range("I665").Value = "b"
$ws.Range("J665").Value = "Acknowledge (Backchannel)"
$ws.Range("I668").Value = "sv"
$ws.Range("J668").Value = "Statement-opinion"
$ws.Range("I673").Value = "sd"
$ws.Range("J673").Value = "Statement-non-opinion"
